$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.528.84'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.148.51'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '531.33'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.67'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.97%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +11.97%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.33'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.62%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.425'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.86%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.109'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.92%  '
$ws.Range("E12").Value = '  +2.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.690.92'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.75'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.92%  '
$ws.Range("E15").Value = '  +4.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.567.52'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.150.09'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.83%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.20'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.99'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.17%  '
$ws.Range("E20").Value = '  +1.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '370.29'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.51%  '
$ws.Range("E22").Value = '  +1.50%  '
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.68'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.515'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.166'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.97'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +10.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0859'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.87'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("E31").Value = '  +0.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.83'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.75%  '
$ws.Range("E33").Value = '  +4.98%  '
$ws.Range("E34").Value = '  +1.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '159.31'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.32%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.28'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.10%  '
$ws.Range("E37").Value = '  +7.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.20'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.656.35'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +10.84%  '
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("E41").Value = '  +1.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.19'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +4.59%  '
$ws.Range("E43").Value = '  +1.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.56'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +4.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0284'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +6.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.191.54'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.92%  '
$ws.Range("E48").Value = '  +11.72%  '
$ws.Range("E49").Value = '  +1.74%  '
$ws.Range("E50").Value = '  +2.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.05'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.18%  '
